$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = -21.107
$ws.Range("A10").Value = -20.945
$ws.Range("A12").Value = -21.694
$ws.Range("D13").Value = -7.831999999999999
$ws.Range("A18").Value = -21.694
